$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the average_doctor / average_doctor_old header labels (BP1 <-> BQ1)
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Update computed statistics cells (rows 4-13) per the new
# "harvard case classification" data pull

# Row 4
$ws.Range("AI4").Value = 0.241
$ws.Range("AJ4").Value = 0.068
$ws.Range("AK4").Value = 0.261
$ws.Range("AU4").Value = 0.157
$ws.Range("AV4").Value = 0.029
$ws.Range("AW4").Value = 0.171
$ws.Range("BA4").Value = 1.925
$ws.Range("BB4").Value = 0.158
$ws.Range("BC4").Value = 0.398
$ws.Range("BG4").Value = 0.72
$ws.Range("BH4").Value = 0.135
$ws.Range("BI4").Value = 0.367
$ws.Range("BM4").Value = 0.6850000000000001
$ws.Range("BN4").Value = 0.08500000000000001
$ws.Range("BO4").Value = 0.292
$ws.Range("BP4").Value = 0.642
$ws.Range("BQ4").Value = 0.646
$ws.Range("E4").Value = 0.409
$ws.Range("F4").Value = 0.074
$ws.Range("G4").Value = 0.272
$ws.Range("N4").Value = 0.411
$ws.Range("O4").Value = 0.064
$ws.Range("P4").Value = 0.253
$ws.Range("W4").Value = 0.236
$ws.Range("X4").Value = 0.108
$ws.Range("Y4").Value = 0.329

# Row 5
$ws.Range("AI5").Value = 0.287
$ws.Range("AJ5").Value = 0.099
$ws.Range("AK5").Value = 0.314
$ws.Range("AU5").Value = 0.314
$ws.Range("AV5").Value = 0.106
$ws.Range("AW5").Value = 0.326
$ws.Range("BA5").Value = 1.349
$ws.Range("BB5").Value = 0.083
$ws.Range("BC5").Value = 0.288
$ws.Range("BG5").Value = 0.402
$ws.Range("BH5").Value = 0.046
$ws.Range("BI5").Value = 0.214
$ws.Range("BM5").Value = 0.5649999999999999
$ws.Range("BN5").Value = 0.073
$ws.Range("BO5").Value = 0.27
$ws.Range("BP5").Value = 0.45
$ws.Range("BQ5").Value = 0.445
$ws.Range("E5").Value = 0.524
$ws.Range("F5").Value = 0.08400000000000001
$ws.Range("G5").Value = 0.289
$ws.Range("N5").Value = 0.742
$ws.Range("O5").Value = 0.08
$ws.Range("P5").Value = 0.283
$ws.Range("W5").Value = 0.238
$ws.Range("X5").Value = 0.113
$ws.Range("Y5").Value = 0.336

# Row 6
$ws.Range("AI6").Value = 0.262
$ws.Range("AU6").Value = 0.209
$ws.Range("BA6").Value = 1.575
$ws.Range("BG6").Value = 0.516
$ws.Range("BM6").Value = 0.619
$ws.Range("BP6").Value = 0.525
$ws.Range("BQ6").Value = 0.524
$ws.Range("E6").Value = 0.459
$ws.Range("N6").Value = 0.529
$ws.Range("W6").Value = 0.237

# Row 7
$ws.Range("AI7").Value = 0.276
$ws.Range("AU7").Value = 0.262
$ws.Range("BA7").Value = 1.43
$ws.Range("BG7").Value = 0.441
$ws.Range("BM7").Value = 0.586
$ws.Range("BP7").Value = 0.477
$ws.Range("BQ7").Value = 0.473
$ws.Range("E7").Value = 0.496
$ws.Range("N7").Value = 0.639
$ws.Range("W7").Value = 0.238

# Row 8
$ws.Range("AI8").Value = 0.263
$ws.Range("AJ8").Value = 0.1
$ws.Range("AK8").Value = 0.316
$ws.Range("AU8").Value = 0.25
$ws.Range("AV8").Value = 0.075
$ws.Range("AW8").Value = 0.275
$ws.Range("BA8").Value = 1.67
$ws.Range("BG8").Value = 0.544
$ws.Range("BH8").Value = 0.097
$ws.Range("BI8").Value = 0.312
$ws.Range("BM8").Value = 0.6830000000000001
$ws.Range("BN8").Value = 0.067
$ws.Range("BO8").Value = 0.26
$ws.Range("BP8").Value = 0.5570000000000001
$ws.Range("BQ8").Value = 0.5679999999999999
$ws.Range("E8").Value = 0.556
$ws.Range("F8").Value = 0.11
$ws.Range("G8").Value = 0.332
$ws.Range("N8").Value = 0.751
$ws.Range("O8").Value = 0.06900000000000001
$ws.Range("P8").Value = 0.262
$ws.Range("W8").Value = 0.247
$ws.Range("X8").Value = 0.116
$ws.Range("Y8").Value = 0.341

# Row 9
$ws.Range("AI9").Value = 0.146
$ws.Range("AJ9").Value = 0.125
$ws.Range("AK9").Value = 0.353
$ws.Range("BA9").Value = 1.561
$ws.Range("BB9").Value = 0.243
$ws.Range("BC9").Value = 0.493
$ws.Range("BG9").Value = 0.5610000000000001
$ws.Range("BH9").Value = 0.246
$ws.Range("BI9").Value = 0.496
$ws.Range("BM9").Value = 0.585
$ws.Range("BN9").Value = 0.243
$ws.Range("BO9").Value = 0.493
$ws.Range("BQ9").Value = 0.52
$ws.Range("E9").Value = 0.463
$ws.Range("F9").Value = 0.249
$ws.Range("G9").Value = 0.499
$ws.Range("N9").Value = 0.634
$ws.Range("O9").Value = 0.232
$ws.Range("P9").Value = 0.482
$ws.Range("W9").Value = 0.146
$ws.Range("X9").Value = 0.125
$ws.Range("Y9").Value = 0.353

# Row 10
$ws.Range("AI10").Value = 0.293
$ws.Range("AJ10").Value = 0.207
$ws.Range("AK10").Value = 0.455
$ws.Range("AU10").Value = 0.244
$ws.Range("AV10").Value = 0.184
$ws.Range("AW10").Value = 0.429
$ws.Range("BA10").Value = 1.927
$ws.Range("BB10").Value = 0.249
$ws.Range("BC10").Value = 0.499
$ws.Range("BG10").Value = 0.61
$ws.Range("BH10").Value = 0.238
$ws.Range("BI10").Value = 0.488
$ws.Range("BM10").Value = 0.854
$ws.Range("BN10").Value = 0.125
$ws.Range("BO10").Value = 0.353
$ws.Range("BP10").Value = 0.642
$ws.Range("BQ10").Value = 0.667
$ws.Range("E10").Value = 0.61
$ws.Range("F10").Value = 0.238
$ws.Range("G10").Value = 0.488
$ws.Range("N10").Value = 0.854
$ws.Range("O10").Value = 0.125
$ws.Range("P10").Value = 0.353
$ws.Range("W10").Value = 0.293
$ws.Range("X10").Value = 0.207
$ws.Range("Y10").Value = 0.455

# Row 11
$ws.Range("AI11").Value = 0.293
$ws.Range("AJ11").Value = 0.207
$ws.Range("AK11").Value = 0.455
$ws.Range("AU11").Value = 0.341
$ws.Range("AV11").Value = 0.225
$ws.Range("AW11").Value = 0.474
$ws.Range("BA11").Value = 1.927
$ws.Range("BB11").Value = 0.249
$ws.Range("BC11").Value = 0.499
$ws.Range("BG11").Value = 0.61
$ws.Range("BH11").Value = 0.238
$ws.Range("BI11").Value = 0.488
$ws.Range("BM11").Value = 0.854
$ws.Range("BN11").Value = 0.125
$ws.Range("BO11").Value = 0.353
$ws.Range("BP11").Value = 0.642
$ws.Range("BQ11").Value = 0.667
$ws.Range("E11").Value = 0.634
$ws.Range("F11").Value = 0.232
$ws.Range("G11").Value = 0.482
$ws.Range("N11").Value = 0.878
$ws.Range("O11").Value = 0.107
$ws.Range("P11").Value = 0.327
$ws.Range("W11").Value = 0.293
$ws.Range("X11").Value = 0.207
$ws.Range("Y11").Value = 0.455

# Row 12
$ws.Range("AI12").Value = 1.917
$ws.Range("AJ12").Value = 0.91
$ws.Range("AK12").Value = 0.954
$ws.Range("AU12").Value = 2.571
$ws.Range("AV12").Value = 1.673
$ws.Range("AW12").Value = 1.294
$ws.Range("BA12").Value = 3.667
$ws.Range("BB12").Value = 0.238
$ws.Range("BC12").Value = 0.488
$ws.Range("BG12").Value = 1.08
$ws.Range("BH12").Value = 0.074
$ws.Range("BI12").Value = 0.271
$ws.Range("BM12").Value = 1.429
$ws.Range("BN12").Value = 0.473
$ws.Range("BO12").Value = 0.6879999999999999
$ws.Range("BP12").Value = 1.222
$ws.Range("BQ12").Value = 1.275
$ws.Range("E12").Value = 1.462
$ws.Range("F12").Value = 0.71
$ws.Range("G12").Value = 0.843
$ws.Range("N12").Value = 1.541
$ws.Range("O12").Value = 1.113
$ws.Range("P12").Value = 1.055
$ws.Range("W12").Value = 1.75
$ws.Range("X12").Value = 0.6879999999999999
$ws.Range("Y12").Value = 0.829

# Row 13
$ws.Range("AI13").Value = 1.365
$ws.Range("AJ13").Value = 0.41
$ws.Range("AK13").Value = 0.64
$ws.Range("AU13").Value = 2.327
$ws.Range("AV13").Value = 0.792
$ws.Range("AW13").Value = 0.89
$ws.Range("BA13").Value = 2.516
$ws.Range("BB13").Value = 0.299
$ws.Range("BC13").Value = 0.547
$ws.Range("BG13").Value = 0.628
$ws.Range("BH13").Value = 0.094
$ws.Range("BI13").Value = 0.307
$ws.Range("BM13").Value = 0.986
$ws.Range("BN13").Value = 0.388
$ws.Range("BO13").Value = 0.623
$ws.Range("BP13").Value = 0.839
$ws.Range("BQ13").Value = 0.768
$ws.Range("E13").Value = 1.652
$ws.Range("F13").Value = 0.865
$ws.Range("G13").Value = 0.93
$ws.Range("N13").Value = 2.205
$ws.Range("O13").Value = 0.928
$ws.Range("P13").Value = 0.963
$ws.Range("W13").Value = 1.078
$ws.Range("X13").Value = 0.181
$ws.Range("Y13").Value = 0.426
